$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.419.31"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.02"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.13"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4738"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2901"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06490"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.93"
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "97.81"
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07707"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7346"
$ws.Range("E13").Value = "  +7.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.870.78"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.112"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.69"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.403.95"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.36"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007547"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.115.24"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.219"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.160"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.262"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.15"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.82"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.924"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1001"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.365"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.508"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.305"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.174"
$ws.Range("E33").Value = "  +4.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04824"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.117"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6959"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01856"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.748"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.293"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.31"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.966"
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4182"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8350"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.64"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.201"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.009"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.29"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "921.43"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05635"
$ws.Range("E51").Value = "  +1.43%  "
